$d = $word.ActiveDocument

# The paragraph reads (in part): "... (one dumbbell for each arm). On flat
# bench press, she can lift 75 lbs. ..." and needs to become "... On
# barbell bench press ...". The word "flat" also shows up elsewhere in the
# same paragraph (e.g. "flat dumbbell press"), so first locate a long,
# unique phrase that pins down the exact occurrence we must change, then
# work out the precise character range of just the word "flat" inside it.
$anchor = $d.Content.Duplicate
$found = $anchor.Find.Execute(
    "On flat bench press, she can lift 75",  # FindText
    $true,                                    # MatchCase
    $false,                                   # MatchWholeWord
    $false,                                   # MatchWildcards
    $false,                                   # MatchSoundsLike
    $false,                                   # MatchAllWordForms
    $true,                                    # Forward
    1,                                         # Wrap (wdFindContinue)
    $false,                                   # Format
    "",                                       # ReplaceWith
    0                                          # Replace (wdReplaceNone)
)

if ($found) {
    # "On " is 3 characters, so "flat" begins right after it and is itself
    # 4 characters long.
    $flatStart = $anchor.Start + 3
    $flatEnd = $flatStart + 4

    $flatRng = $d.Range($flatStart, $flatEnd)

    if ($flatRng.Text -eq "flat") {
        # Swap the word itself; the rest of the sentence is untouched.
        $flatRng.Text = "barbell"

        # The newly typed word would otherwise silently merge back into the
        # surrounding run because its formatting is identical. Toggle a
        # character property on just the new word to force it to live in
        # its own run, split off from the text before and after it.
        $newWordRng = $d.Range($flatStart, $flatStart + 7)
        $newWordRng.Font.Bold = 1
        $newWordRng.Font.Bold = 0
    }
}
